$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Avion3"
$ws.Range("B4").Value = "v1000"

$ws.Range("A5").Value = "Version RPB"
$ws.Range("B5").Value = "RPB"
